$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Add the new "Enterprise" worksheet at the end of the tab strip (after
# "Paint", the last existing sheet) so it becomes the 4th / rightmost
# sheet and the newly active one.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Enterprise"

# ---------------------------------------------------------------------
# Header row - same column headers used by the other shop sheets
# (Arbeit / Realty / Paint): code, name, price, per_second, buy_status,
# profit.
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "code"
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "price"
$ws.Range("D1").Value = "per_second"
$ws.Range("E1").Value = "buy_status"
$ws.Range("F1").Value = "profit"

# ---------------------------------------------------------------------
# Data rows: 15 new "enterprise" buttons.
# ---------------------------------------------------------------------
$names = @("테쓸라","bwm","벤츄","폭소바겐","스타박스","도요토","넷플리즈","막도널드","카카콜라","탄센트","페이스쿡","그글","아마도존","마이크로하드","애플들")
$prices = @(10000000000000,20000000000000,30000000000000,50000000000000,100000000000000,200000000000000,300000000000000,500000000000000,1000000000000000,2000000000000000,3000000000000000,5000000000000000,10000000000000000,20000000000000000,30000000000000000)
$perSecond = @(100000000,200000000,300000000,500000000,100000000,2000000000,3000000000,5000000000,10000000000,20000000000,30000000000,50000000000,100000000000,200000000000,300000000000)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $i + 1
    $ws.Cells.Item($row, 2).Value = $names[$i]
    $ws.Cells.Item($row, 3).Value = $prices[$i]
    $ws.Cells.Item($row, 4).Value = $perSecond[$i]
    $ws.Cells.Item($row, 5).Value = $false
    $ws.Cells.Item($row, 6).Value = 0
}

# Column widths (bestFit-style, matching the other shop sheets).
$ws.Columns.Item(3).ColumnWidth = 8.3
$ws.Columns.Item(4).ColumnWidth = 11.8

# Page setup to match the other data sheets.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Selection left on the new sheet, matching the authored workbook.
$ws.Range("F3:F16").Select() | Out-Null
